$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.8203074518761176
$ws.Range("J2").Value = 0.8725723693674974
$ws.Range("M2").Value = 30.46625333333334
$ws.Range("N2").Value = 91.39876000000001
$ws.Range("O2").Value = 0.2185380492512374
$ws.Range("P2").Value = 0.2331534018544084
$ws.Range("Q2").Value = 1.663122303213334
$ws.Range("R2").Value = 14.96810072892
$ws.Range("S2").Value = 0.1792683903192601
$ws.Range("T2").Value = 0.2034432162821934

$ws.Range("I3").Value = 0.8203074518761176
$ws.Range("J3").Value = 0.8725723693674974
$ws.Range("O3").Value = 0.2491807703757967
$ws.Range("P3").Value = 0.2658454419670822
$ws.Range("S3").Value = 0.2044048428034977
$ws.Range("T3").Value = 0.2319693871827664

$ws.Range("I4").Value = 0.8203074518761176
$ws.Range("J4").Value = 0.8725723693674974
$ws.Range("M4").Value = 23.69037333333334
$ws.Range("N4").Value = 71.07112000000001
$ws.Range("O4").Value = 0.1699338582153697
$ws.Range("P4").Value = 0.181298667526812
$ws.Range("Q4").Value = 1.293233789893333
$ws.Range("R4").Value = 11.63910410904
$ws.Range("S4").Value = 0.1393980102201274
$ws.Range("T4").Value = 0.1581962078870405

$ws.Range("I5").Value = 0.8203074518761176
$ws.Range("J5").Value = 0.8725723693674974
$ws.Range("M5").Value = 26.2168665
$ws.Range("N5").Value = 52.433733
$ws.Range("O5").Value = 0.18805669340777
$ws.Range("P5").Value = 0.1337556791894743
$ws.Range("Q5").Value = 1.4311525253685
$ws.Range("R5").Value = 8.586915152211001
$ws.Range("S5").Value = 0.1542643069775761
$ws.Range("T5").Value = 0.1167115099067185

$ws.Range("I6").Value = 0.8203074518761176
$ws.Range("J6").Value = 0.8725723693674974
$ws.Range("M6").Value = 24.297748
$ws.Range("N6").Value = 72.893244
$ws.Range("O6").Value = 0.1742906287498262
$ws.Range("P6").Value = 0.1859468094622229
$ws.Range("Q6").Value = 1.326389765572
$ws.Range("R6").Value = 11.937507890148
$ws.Range("S6").Value = 0.1429719015556563
$ws.Range("T6").Value = 0.1622520481087785

$ws.Range("G7").Value = 0.011958
$ws.Range("H7").Value = 0.023916
$ws.Range("I7").Value = 0.1796925481238824
$ws.Range("J7").Value = 0.1274276306325027
$ws.Range("M7").Value = 30.46625333333334
$ws.Range("N7").Value = 91.39876000000001
$ws.Range("O7").Value = 0.2185380492512374
$ws.Range("P7").Value = 0.2331534018544084
$ws.Range("Q7").Value = 0.36431545736
$ws.Range("R7").Value = 2.18589274416
$ws.Range("S7").Value = 0.03926965893197736
$ws.Range("T7").Value = 0.02971018557221502

$ws.Range("G8").Value = 0.011958
$ws.Range("H8").Value = 0.023916
$ws.Range("I8").Value = 0.1796925481238824
$ws.Range("J8").Value = 0.1274276306325027
$ws.Range("O8").Value = 0.2491807703757967
$ws.Range("P8").Value = 0.2658454419670822
$ws.Range("Q8").Value = 0.4153986302879999
$ws.Range("R8").Value = 2.492391781728
$ws.Range("S8").Value = 0.04477592757229893
$ws.Range("T8").Value = 0.03387605478431577

$ws.Range("G9").Value = 0.011958
$ws.Range("H9").Value = 0.023916
$ws.Range("I9").Value = 0.1796925481238824
$ws.Range("J9").Value = 0.1274276306325027
$ws.Range("M9").Value = 23.69037333333334
$ws.Range("N9").Value = 71.07112000000001
$ws.Range("O9").Value = 0.1699338582153697
$ws.Range("P9").Value = 0.181298667526812
$ws.Range("Q9").Value = 0.2832894843200001
$ws.Range("R9").Value = 1.69973690592
$ws.Range("S9").Value = 0.03053584799524233
$ws.Range("T9").Value = 0.02310245963977151

$ws.Range("G10").Value = 0.011958
$ws.Range("H10").Value = 0.023916
$ws.Range("I10").Value = 0.1796925481238824
$ws.Range("J10").Value = 0.1274276306325027
$ws.Range("M10").Value = 26.2168665
$ws.Range("N10").Value = 52.433733
$ws.Range("O10").Value = 0.18805669340777
$ws.Range("P10").Value = 0.1337556791894743
$ws.Range("Q10").Value = 0.313501289607
$ws.Range("R10").Value = 1.254005158428
$ws.Range("S10").Value = 0.0337923864301939
$ws.Range("T10").Value = 0.01704416928275586

$ws.Range("G11").Value = 0.011958
$ws.Range("H11").Value = 0.023916
$ws.Range("I11").Value = 0.1796925481238824
$ws.Range("J11").Value = 0.1274276306325027
$ws.Range("M11").Value = 24.297748
$ws.Range("N11").Value = 72.893244
$ws.Range("O11").Value = 0.1742906287498262
$ws.Range("P11").Value = 0.1859468094622229
$ws.Range("Q11").Value = 0.290552470584
$ws.Range("R11").Value = 1.743314823504
$ws.Range("S11").Value = 0.03131872719416987
$ws.Range("T11").Value = 0.0236947613534445
